# Mise à jour de l'application
# Adds 7 new training-session rows (J-2, 2025-09-18) to the data sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row data: Type, Date, Periode, MD, Nom, Poste, Temps joue, H..V (15 numeric stats)
$rows = @(
  @{ E="Kamal Bafounta";    F="center midfield"; G="00:46:52"; H=2.73; I=0.33; J=2.39;  K=0.13; L=0.12; M=0.07; N=0;    O=5; P=3.37; Q=29.42; R=4.2;  S=20; T=1; U=3;  V=5 },
  @{ E="Naim Ighbane";      F="center back";     G="00:32:22"; H=2.39; I=0.26; J=2.13;  K=0.15; L=0.11; M=0;    N=0;    O=1; P=3.82; Q=25.69; R=4.74; S=17; T=3; U=4;  V=1 },
  @{ E="Yanis Berrached";   F="center midfield"; G="00:30:10"; H=1.86; I=0.25; J=1.6;   K=0.13; L=0.1;  M=0.03; N=0;    O=3; P=3.59; Q=28.02; R=4.31; S=10; T=1; U=2;  V=0 },
  @{ E="Malik Boussaid";    F="right back";      G="00:43:45"; H=3.01; I=0.49; J=2.51;  K=0.24; L=0.16; M=0.1;  N=0;    O=7; P=4.03; Q=29.63; R=4.15; S=16; T=2; U=10; V=1 },
  @{ E="Omar Benyounes";    F="center midfield"; G="00:47:14"; H=3;    I=0.33; J=2.67;  K=0.11; L=0.13; M=0.09; N=0;    O=5; P=3.76; Q=29.37; R=3.99; S=18; T=0; U=7;  V=0 },
  @{ E="Hedi Nasri";        F="right back";      G="00:32:38"; H=2.21; I=0.3;  J=1.91;  K=0.07; L=0.09; M=0.12; N=0.01; O=8; P=3.96; Q=30.73; R=4.27; S=4;  T=1; U=3;  V=0 },
  @{ E="Mattheo Haon";      F="right back";      G="00:45:34"; H=2.6;  I=0.21; J=2.38;  K=0.13; L=0.08; M=0.01; N=0;    O=2; P=3.29; Q=26.64; R=4.42; S=11; T=2; U=4;  V=1 }
)

$startRow = 515
$sourceRow = 514

for ($i = 0; $i -lt $rows.Count; $i++) {
    $r = $startRow + $i
    $data = $rows[$i]

    $ws.Cells.Item($r, 1).Value = "Entrainement"
    $ws.Cells.Item($r, 2).Value = 45918
    $ws.Cells.Item($r, 3).Value = "Global"
    $ws.Cells.Item($r, 4).Value = "J-2"
    $ws.Cells.Item($r, 5).Value = $data.E
    $ws.Cells.Item($r, 6).Value = $data.F
    $ws.Cells.Item($r, 7).Value = $data.G
    $ws.Cells.Item($r, 8).Value = $data.H
    $ws.Cells.Item($r, 9).Value = $data.I
    $ws.Cells.Item($r, 10).Value = $data.J
    $ws.Cells.Item($r, 11).Value = $data.K
    $ws.Cells.Item($r, 12).Value = $data.L
    $ws.Cells.Item($r, 13).Value = $data.M
    $ws.Cells.Item($r, 14).Value = $data.N
    $ws.Cells.Item($r, 15).Value = $data.O
    $ws.Cells.Item($r, 16).Value = $data.P
    $ws.Cells.Item($r, 17).Value = $data.Q
    $ws.Cells.Item($r, 18).Value = $data.R
    $ws.Cells.Item($r, 19).Value = $data.S
    $ws.Cells.Item($r, 20).Value = $data.T
    $ws.Cells.Item($r, 21).Value = $data.U
    $ws.Cells.Item($r, 22).Value = $data.V

    # Copy number formats (styles) from the last existing row so we reuse
    # the same style indexes instead of generating new ones.
    $ws.Cells.Item($sourceRow, 2).Copy()
    $ws.Cells.Item($r, 2).PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)

    $ws.Cells.Item($sourceRow, 4).Copy()
    $ws.Cells.Item($r, 4).PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)
}

$excel.CutCopyMode = 0

# Update the selection to match where the user ended up after entering data.
$ws.Range("E524").Select()

$wb.Save()
